# Updates the cryptos price table (rows 2-51) with the latest scraped
# values from the GitHub Actions refresh job. Coin/Link/Price/Volume(1h)
# cells are plain text in this sheet; a leading apostrophe is used on
# numeric-looking Price values so Excel keeps storing them as text
# (e.g. "1.000") instead of silently re-parsing them as numbers and
# dropping the trailing zero / formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.181.36'
$ws.Range("E2").Value = '  -3.20%  '

$ws.Range("D3").Value = '1.606.57'
$ws.Range("E3").Value = '  -2.68%  '

$ws.Range("D4").Value = '''1.000'

$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").Value = '''302.30'
$ws.Range("E6").Value = '  -2.09%  '

$ws.Range("E7").Value = '  -2.64%  '

$ws.Range("D8").Value = '''0.3667'
$ws.Range("E8").Value = '  -4.12%  '

$ws.Range("D9").Value = '''49.02'
$ws.Range("E9").Value = '  -5.35%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '''1.270'
$ws.Range("E10").Value = '  -5.98%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.08146'
$ws.Range("E11").Value = '  -3.63%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '''1.000'
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("D13").Value = '''23.18'
$ws.Range("E13").Value = '  -2.84%  '

$ws.Range("D14").Value = '''6.630'
$ws.Range("E14").Value = '  -6.45%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.00001262'
$ws.Range("E15").Value = '  -3.83%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''7.407'
$ws.Range("E16").Value = '  -6.96%  '

$ws.Range("D17").Value = '1.611.30'
$ws.Range("E17").Value = '  -2.29%  '

$ws.Range("D18").Value = '''91.69'
$ws.Range("E18").Value = '  -3.21%  '

$ws.Range("D19").Value = '''0.06817'
$ws.Range("E19").Value = '  -2.01%  '

$ws.Range("D20").Value = '''18.33'
$ws.Range("E20").Value = '  -6.94%  '

$ws.Range("D21").Value = '''6.580'
$ws.Range("E21").Value = '  -5.76%  '

$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = '''13.03'
$ws.Range("E23").Value = '  -5.63%  '

$ws.Range("D24").Value = '23.205.78'
$ws.Range("E24").Value = '  -3.15%  '

$ws.Range("D25").Value = '''2.351'
$ws.Range("E25").Value = '  -4.07%  '

$ws.Range("D26").Value = '''2.865'
$ws.Range("E26").Value = '  -5.30%  '

$ws.Range("D27").Value = '''21.14'
$ws.Range("E27").Value = '  -4.57%  '

$ws.Range("D28").Value = '''151.39'
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("D29").Value = '''5.241'
$ws.Range("E29").Value = '  -2.68%  '

$ws.Range("D30").Value = '''132.38'
$ws.Range("E30").Value = '  -4.60%  '

$ws.Range("D31").Value = '''2.443'
$ws.Range("E31").Value = '  -3.40%  '

$ws.Range("D32").Value = '''6.926'
$ws.Range("E32").Value = '  -11.49%  '

$ws.Range("D33").Value = '1.787.69'
$ws.Range("E33").Value = '  -2.32%  '

$ws.Range("D34").Value = '''0.9641'
$ws.Range("E34").Value = '  -7.38%  '

$ws.Range("D35").Value = '''0.07783'
$ws.Range("E35").Value = '  -3.52%  '

$ws.Range("D36").Value = '''0.02787'
$ws.Range("E36").Value = '  -5.60%  '

$ws.Range("D37").Value = '''6.299'
$ws.Range("E37").Value = '  -5.45%  '

$ws.Range("D38").Value = '''0.2558'
$ws.Range("E38").Value = '  -4.18%  '

$ws.Range("D39").Value = '''10.19'
$ws.Range("E39").Value = '  -5.46%  '

$ws.Range("D40").Value = '''0.08930'
$ws.Range("E40").Value = '  -1.84%  '

$ws.Range("D41").Value = '''1.391'
$ws.Range("E41").Value = '  -2.19%  '

$ws.Range("D42").Value = '''0.7207'
$ws.Range("E42").Value = '  -5.09%  '

$ws.Range("D43").Value = '''12.82'
$ws.Range("E43").Value = '  -4.31%  '

$ws.Range("D44").Value = '''15.57'
$ws.Range("E44").Value = '  -4.76%  '

$ws.Range("D45").Value = '''0.6694'
$ws.Range("E45").Value = '  -4.50%  '

$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '''0.9997'
$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.311'
$ws.Range("E47").Value = '  -6.30%  '

$ws.Range("D48").Value = '''3.985'
$ws.Range("E48").Value = '  -2.34%  '

$ws.Range("D49").Value = '''0.07989'
$ws.Range("E49").Value = '  -3.95%  '

$ws.Range("D50").Value = '''131.69'
$ws.Range("E50").Value = '  -2.54%  '

$ws.Range("D51").Value = '''1.212'
$ws.Range("E51").Value = '  +0.21%  '
